# SectorGroup.xlsx fix-up:
# The columns were mislabeled/shifted: E (group-code), F (group-name) and
# G (category-name) actually held each other's data one step out of place.
# For every row (including the header) the correct values are obtained by
# rotating the E/F/G values one step to the right:
#   new E = old G
#   new F = old E
#   new G = old F

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item(1, 1).End(-4121).Row   # -4121 = xlDown

for ($r = 1; $r -le $lastRow; $r++) {
    $eOld = $ws.Cells.Item($r, 5).Value2
    $fOld = $ws.Cells.Item($r, 6).Value2
    $gOld = $ws.Cells.Item($r, 7).Value2

    $ws.Cells.Item($r, 5).Value2 = $gOld
    $ws.Cells.Item($r, 6).Value2 = $eOld
    $ws.Cells.Item($r, 7).Value2 = $fOld
}
